# Grow the dark "code block" panels on slides 5, 6 and 14 (they were
# clipping their Courier New text) and push the shapes underneath them
# down by the same amount so nothing overlaps.

$p = $ppt.ActivePresentation

# --- Slide 5 ("Code Block Support": Python + JavaScript examples) ---
$s5 = $p.Slides.Item(5)

# Python code block grows taller (top stays put)
$s5.Shapes.Item(3).Height = 169.5

# "JavaScript Example" label shifts down to follow
$s5.Shapes.Item(4).Top = 250.5

# JavaScript code block shifts down and also grows taller
$s5.Shapes.Item(5).Top = 276.75
$s5.Shapes.Item(5).Height = 169.5

# "SQL Example" label shifts down to follow
$s5.Shapes.Item(6).Top = 457.49993

# --- Slide 6 (SQL code block) ---
$s6 = $p.Slides.Item(6)

# SQL code block grows taller (top stays put)
$s6.Shapes.Item(1).Height = 201.74992125984252

# --- Slide 14 ("Technical Details": API usage) ---
$s14 = $p.Slides.Item(14)

# API usage code block grows taller (top stays put)
$s14.Shapes.Item(4).Height = 169.5

# "Configuration Options" label shifts down to follow
$s14.Shapes.Item(5).Top = 273.0

# Configuration bullet list shifts down to follow
$s14.Shapes.Item(6).Top = 299.25

# "End of demonstration" line shifts down to follow
$s14.Shapes.Item(7).Top = 360.0
